$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing "key column composition" entry in row 13
# (was "Потребность-Дата отгрузки-Дата прихода", now "Потребность-Уведомление")
$ws.Range("C13").Value = "Потребность-Уведомление"

# Add a new row 14 with the new key-column attribute mapping
$ws.Range("A14").Value = "6532ec8d-56cd-ed11-9165-005056b6948b"
$ws.Range("B14").Value = 2
$ws.Range("C14").Value = "Резервирование/Планирование закупок/Спецификация/РасходныйДокумент/Поступление/Уведомление об отгрузке.Номер"

# Match the resulting cell selection
$ws.Range("C10").Select()
